$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Captured_Values")

$text = "Real Programmers Count 0123456789 From Zero"

for ($r = 21; $r -le 35; $r++) {
    $ws.Cells.Item($r, 1).Value = 123456789
    $ws.Cells.Item($r, 2).Value = $text
}
